$d = $word.ActiveDocument

# --- "Towers > 5 Types" block: strip the detailed stat call-outs back down
# to a plain template list (Cannon / Lasers / Electricity / Sniper / Rockets),
# matching the target diff.
#
# Original (1-indexed) paragraphs 34-47:
#  34 "5 Types: (Can attack, level 1) (Can attack, level 2) (Can attack, level 3)"
#  35 "Cannon (Basic, Basic, Basic + Cloaked)"
#  36   "Strong against: Light, Medium"                         (ilvl3, has _GoBack bookmark)
#  37   "Weak against: Heavy"                                   (ilvl3)
#  38 "Lasers (Basic + Cloaked, Basic + Cloaked, Basic + Cloaked)"
#  39   "Strong against: Heavy"                                 (ilvl3)
#  40   "Weak against: Light"                                   (ilvl3)
#  41 "Electricity (Basic, Basic, Basic)"
#  42   "Strong against: Heavy, Medium, Light, Close enemies (for chain)" (ilvl3)
#  43 "Sniper (Basic + Cloaked, Basic + Cloaked, Basic + Cloaked)"
#  44   "Strong against: Heavy, Invisible"                      (ilvl3)
#  45   "Weak against: Medium, Light"                           (ilvl3)
#  46 "Rockets (Basic, Basic, Basic + Cloaked)"
#  47   "Strong against: Clustered enemies (for area of effect damage), Light, Medium, Heavy" (ilvl3)
#
# Work from the bottom of the block upward so earlier paragraph indices stay valid
# while we delete paragraphs.

# 47: delete "Strong against: Clustered..." sub-bullet (under Rockets)
$d.Paragraphs(47).Range.Delete()

# 46: "Rockets (Basic, Basic, Basic + Cloaked)" -> "Rockets"
$r = $d.Paragraphs(46).Range.Duplicate
$r.Find.Execute(" (Basic, Basic, Basic + Cloaked)")
$r.Delete()

# 45: delete "Weak against: Medium, Light"
$d.Paragraphs(45).Range.Delete()
# 44: delete "Strong against: Heavy, Invisible"
$d.Paragraphs(44).Range.Delete()

# 43: "Sniper (Basic + Cloaked, Basic + Cloaked, Basic + Cloaked)" -> "Sniper"
$r = $d.Paragraphs(43).Range.Duplicate
$r.Find.Execute(" (Basic + Cloaked, Basic + Cloaked, Basic + Cloaked)")
$r.Delete()
# re-anchor the _GoBack bookmark right after "Sniper" (matches the target XML)
$bm = $d.Content
$bm.Find.Execute("Sniper")
$bm.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bm)

# 42: delete "Strong against: Heavy, Medium, Light, Close enemies (for chain)"
$d.Paragraphs(42).Range.Delete()

# 41: "Electricity (Basic, Basic, Basic)" -> "Electricity"
$r = $d.Paragraphs(41).Range.Duplicate
$r.Find.Execute(" (Basic, Basic, Basic)")
$r.Delete()

# 40: delete "Weak against: Light"
$d.Paragraphs(40).Range.Delete()
# 39: delete "Strong against: Heavy"
$d.Paragraphs(39).Range.Delete()

# 38: "Lasers (Basic + Cloaked, Basic + Cloaked, Basic + Cloaked)" -> "Lasers"
$r = $d.Paragraphs(38).Range.Duplicate
$r.Find.Execute(" (Basic + Cloaked, Basic + Cloaked, Basic + Cloaked)")
$r.Delete()

# 37: delete "Weak against: Heavy"
$d.Paragraphs(37).Range.Delete()
# 36: delete "Strong against: Light, Medium" (original _GoBack location; removed along with paragraph)
$d.Paragraphs(36).Range.Delete()

# 35: "Cannon (Basic, Basic, Basic + Cloaked)" -> "Cannon"
$r = $d.Paragraphs(35).Range.Duplicate
$r.Find.Execute(" (Basic, Basic, Basic + Cloaked)")
$r.Delete()

# 34: "5 Types: (Can attack, level 1) (Can attack, level 2) (Can attack, level 3)"
#     -> "5 Types: " (trailing single space)
$r = $d.Paragraphs(34).Range.Duplicate
$r.Find.Execute(" (Can attack, level 1) (Can attack, level 2) (Can attack, level 3)", $true, $false, $false, $false, $false, $true, 1, $false, " ", 2)
